$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.928.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.638.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'579.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'157.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.76%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.85%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'5.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.98%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.52%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.04%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'28.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.55%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.114.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.23%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'63.847.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.655.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.02%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.14%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'344.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.73%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.11%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'68.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.19%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +9.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +5.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +5.08%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.48%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'584.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.35%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.18%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.40%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.60%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.68%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'19.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'154.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.59%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +9.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D43").Value = "'162.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.35%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'24.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.0591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.31%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.57%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0250"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +2.15%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.794"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.18%  "
$ws.Range("E51").Style = "Normal"
